$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells: "_old" suffix -> "_FV2310", "_new" suffix -> "_FV2404"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = $cell.Text
    if ($text -like "*_old") {
        $cell.Value = ($text -replace "_old$", "_FV2310")
    } elseif ($text -like "*_new") {
        $cell.Value = ($text -replace "_new$", "_FV2404")
    }
}

# 2. Turn the used range into an Excel Table ("Table1") with an autofilter
$range = $ws.Range("A1:U64")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# 3. Freeze the header row (split below row 1)
$ws.Activate() | Out-Null
$ws.Cells.Item(2, 1).Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
